$wb = $excel.ActiveWorkbook

# Sheets (by position): 1 = Mflag, 2 = Pflag, 3 = Sflag, 4 = Uflag
$wsM = $wb.Worksheets.Item(1)
$wsP = $wb.Worksheets.Item(2)
$wsS = $wb.Worksheets.Item(3)
$wsU = $wb.Worksheets.Item(4)

# --- Mflag: drop the stale "/configuration/clientCode/assays" test rows ---
$wsM.Range("A13").ClearContents()
$wsM.Range("A18").ClearContents()
$wsM.Range("A22").ClearContents()

# --- Pflag: point the test cases at the new endpoints ---
$wsP.Range("A3").Value = "/donationInfo/donations"
$wsP.Range("A8").Value = "/testInfo/orderTests"

# --- Sflag: point the test cases at the new endpoints ---
$wsS.Range("A3").Value = "/donationInfo/donations"
$wsS.Range("A8").Value = "/testInfo/orderTests"

# --- Uflag: point the test cases at the new endpoints ---
$wsU.Range("A3").Value = "/donationInfo/donations"
$wsU.Range("A8").Value = "/testInfo/orderTests"

# --- Update each sheet's selection / active cell ---
$wsM.Activate()
$wsM.Range("A3").Select()

$wsP.Activate()
$wsP.Range("A3").Select()

$wsS.Activate()
$wsS.Range("A3").Select()

# Uflag ends up as the active tab, selection on D19
$wsU.Activate()
$wsU.Range("D19").Select()
